$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -12.06679999999999
    6   = -11.4844
    7   = -12.8417
    8   = -12.5971
    16  = -14.39549999999999
    20  = -12.2744
    21  = -13.3515
    28  = -13.1479
    29  = -11.368
    30  = -12.2689
    32  = -11.762
    40  = -12.56660000000001
    46  = -14.35689999999999
    51  = -11.5947
    52  = -11.01230000000001
    57  = -14.50519999999999
    59  = -12.5984
    62  = -14.25289999999999
    66  = -11.2022
    73  = -11.9366
    74  = -12.0079
    77  = -12.57600000000001
    92  = -10.6786
    100 = -11.42909999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
